$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# 1) "Google.Books API" (Product Scope Definition bullet list) -> split run,
#    wrap "Google.Books" in spellStart/spellEnd proofErr marks.
$rng = $d.Content
$rng.Find.Execute("Google.Books API", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$frag1 = '<w:p ' + $wNs + ' w14:paraId="15909EEB" w14:textId="77777777" w:rsidR="003E6BDD" w:rsidRDefault="003E6BDD" w:rsidP="007432EE">' + `
  '<w:pPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Google.Books</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> API</w:t></w:r>' + `
  '</w:p>'
$rng.InsertXML($frag1) | Out-Null

# 2) "Funtionality required:" -> fix spelling to "Functionality required:", split into two runs
#    (no proofErr since the corrected word is not flagged).
$rng = $d.Content
$rng.Find.Execute("Funtionality required:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$frag2 = '<w:p ' + $wNs + ' w14:paraId="731A023B" w14:textId="4185E9C1" w:rsidR="00BD49B5" w:rsidRDefault="00177B00" w:rsidP="007432EE">' + `
  '<w:pPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Functionality</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> required:</w:t></w:r>' + `
  '</w:p>'
$rng.InsertXML($frag2) | Out-Null

# 3) "Postman API calls to Google.Books" -> split run, wrap "Google.Books" in proofErr marks.
$rng = $d.Content
$rng.Find.Execute("Postman API calls to Google.Books", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$frag3 = '<w:p ' + $wNs + ' w14:paraId="78DBDE02" w14:textId="57CD6308" w:rsidR="00D20C7F" w:rsidRDefault="00D20C7F" w:rsidP="00F7616D">' + `
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">Postman API calls to </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Google.Books</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>'
$rng.InsertXML($frag3) | Out-Null

# 4) "Use Google.Books API to generate data in our database" -> split into 3 runs,
#    wrap "Google.Books" in proofErr marks.
$rng = $d.Content
$rng.Find.Execute("Use Google.Books API to generate data in our database", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$frag4 = '<w:p ' + $wNs + ' w14:paraId="50320214" w14:textId="58F07F17" w:rsidR="00E627F6" w:rsidRDefault="00853227" w:rsidP="00E627F6">' + `
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">Use </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Google.Books</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> API to generate data in our database</w:t></w:r>' + `
  '</w:p>'
$rng.InsertXML($frag4) | Out-Null

Write-Host "Applied 4 edits"
